$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ur = $ws.UsedRange
$rowCount = $ur.Rows.Count
$colCount = $ur.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $orig = $cell.Text
        if ($orig -like '*D51*' -or $orig -like '*D64*' -or $orig -like '*D80*' -or $orig -like '*S30*') {
            $new = $orig -replace 'D51', 'D55'
            $new = $new -replace 'D64', 'D69'
            $new = $new -replace 'D80', 'D86'
            $new = $new -replace 'S30', 'S31'
            $cell.Value = $new
        }
    }
}
